$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B:C").Delete()
$ws.Rows("25:32").Delete()
$dates = @(45728, 45729, 45730, 45733, 45734, 45735, 45736, 45737, 45740, 45741, 45742, 45743, 45744, 45748, 45749, 45750, 45751, 45754, 45755, 45756, 45757, 45758, 45761)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}
$ws.Range("A2:A24").NumberFormat = "YYYY-MM-DD"
Write-Output "done"
